$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.012.20'
$ws.Range('E2').Value = '  +0.95%  '

$ws.Range('D3').Value = '1.845.71'
$ws.Range('E3').Value = '  +2.00%  '

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.18%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.620'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.03%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.98%  '

$ws.Range('E9').Value = '  +1.22%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0693'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.83%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0983'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.81%  '

$ws.Range('D12').Value = '2.114.45'
$ws.Range('E12').Value = '  +2.08%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.853.45'
$ws.Range('E13').Value = '  +2.25%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.57%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.674'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.99%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.68'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.93%  '

$ws.Range('D17').Value = '35.008.68'
$ws.Range('E17').Value = '  +0.96%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.91'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.72%  '

$ws.Range('D19').Value = '0.0₃0791'
$ws.Range('E19').Value = '  +0.98%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.08'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.88%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.91%  '

$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('E24').Value = '  +2.71%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.27%  '

$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +20.85%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.09%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.17%  '

$ws.Range('E29').Value = '  +2.57%  '

$ws.Range('E30').Value = '  +0.18%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0554'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.01%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.02%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.97'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.85%  '

$ws.Range('E34').Value = '  +24.87%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.97'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +10.21%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.757'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.77%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.22'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.55%  '

$ws.Range('E38').Value = '  +11.32%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '90.81'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.58%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0200'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.36%  '

$ws.Range('D41').Value = '1.345.87'
$ws.Range('E41').Value = '  +1.33%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.67'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.77%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.46%  '

$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.39'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.19%  '

$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.92%  '

$ws.Range('B46').Value = 'Gas'
$ws.Range('C46').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.35'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +81.40%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0532'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.65%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.35'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.68%  '

$ws.Range('D49').Value = '2.025.92'
$ws.Range('E49').Value = '  +1.37%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.44'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +15.91%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0675'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.81%  '
